# Edit the Products worksheet:
#  - Row 10 (Caymus Cabernet Sauvignon) and Row 12 (Corona Extra) get their
#    name fields tweaked, several boolean flag columns switch from native
#    Excel booleans to literal text "true"/"false", the alcohol % column
#    becomes a plain number, size/volume text changes, the "isTrending"
#    flag becomes literal text "true", and new image_url / date_modified
#    values are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 10: Caymus Cabernet Sauvignon ----------------------------------
$ws.Range("C10").Value = "Caymus Cabernet Sauvignon2"
$ws.Range("D10").Value = "Caymus Cabernet Sauvignon2"

# A leading apostrophe forces these look-alike words to be stored as literal
# text instead of being auto-coerced into the native Excel Boolean type.
$ws.Range("W10").Value = "'false"
$ws.Range("X10").Value = "'false"
$ws.Range("Y10").Value = "'false"
$ws.Range("Z10").Value = "'true"

$ws.Range("AA10").Value = 14.6

$ws.Range("AB10").Value = "75ml"
$ws.Range("AC10").Value = "75ML"

$ws.Range("AS10").Value = "'true"
$ws.Range("AT10").Value = "https://ext.same-assets.com/1701767421/4224202088.png"
$ws.Range("AU10").Value = "2025-03-30T18:18:14.916Z"

# ---- Row 12: Corona Extra -------------------------------------------------
$ws.Range("C12").Value = "Corona Extra45"
$ws.Range("D12").Value = "Corona Extra45"

$ws.Range("W12").Value = "'false"
$ws.Range("X12").Value = "'false"
$ws.Range("Y12").Value = "'false"
$ws.Range("Z12").Value = "'false"

$ws.Range("AA12").Value = 4.5

$ws.Range("AB12").Value = "633ml"
$ws.Range("AC12").Value = "633ML"

$ws.Range("AS12").Value = "'true"
$ws.Range("AT12").Value = "https://ext.same-assets.com/1701767421/2183896642.jpeg"
$ws.Range("AU12").Value = "2025-03-30T18:17:12.776Z"
